$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FX rates")

# Row 13 - USD/IDR
$ws.Range("J13").Value = 14430
$ws.Range("K13").Value = 14360
$ws.Range("L13").Value = 14290
$ws.Range("N13").Value = 14210
$ws.Range("O13").Value = 14140
$ws.Range("P13").Value = 14070

# Row 14 - USD/MYR
$ws.Range("I14").Value = 4.1100000000000003
$ws.Range("J14").Value = 4.09
$ws.Range("K14").Value = 4.08
$ws.Range("L14").Value = 4.0599999999999996

# Row 16 - USD/SGD
$ws.Range("I16").Value = 1.35
$ws.Range("J16").Value = 1.34

# Row 17 - USD/KRW
$ws.Range("I17").Value = 1115
$ws.Range("J17").Value = 1110
$ws.Range("K17").Value = 1100
$ws.Range("L17").Value = 1090
$ws.Range("N17").Value = 1085
$ws.Range("O17").Value = 1075
$ws.Range("P17").Value = 1070
$ws.Range("Q17").Value = 1060

# Row 19 - USD/VND
$ws.Range("I19").Value = 23160
$ws.Range("J19").Value = 23100
$ws.Range("K19").Value = 23050
$ws.Range("L19").Value = 22990
$ws.Range("N19").Value = 22940
$ws.Range("O19").Value = 22890
$ws.Range("P19").Value = 22840
$ws.Range("Q19").Value = 22780

# Row 21 - AUD/USD
$ws.Range("K21").Value = 0.74
$ws.Range("L21").Value = 0.75
$ws.Range("N21").Value = 0.75
$ws.Range("O21").Value = 0.76
$ws.Range("P21").Value = 0.77
$ws.Range("Q21").Value = 0.77

# Row 22 - EUR/USD
$ws.Range("I22").Value = 1.18
$ws.Range("J22").Value = 1.19
$ws.Range("K22").Value = 1.2
$ws.Range("L22").Value = 1.21
$ws.Range("N22").Value = 1.22
$ws.Range("O22").Value = 1.23

# Row 24 - GBP/USD
$ws.Range("I24").Value = 1.3
$ws.Range("J24").Value = 1.31
$ws.Range("K24").Value = 1.32
$ws.Range("L24").Value = 1.33
$ws.Range("N24").Value = 1.34
$ws.Range("O24").Value = 1.35
$ws.Range("P24").Value = 1.36
$ws.Range("Q24").Value = 1.37

# Update the selected cell on the FX rates sheet
$ws.Range("R24").Select()

$excel.CalculateFullRebuild()
